$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 4, pushing the old row 4
# (red, 0, 100, high) down to row 6, and shifting the "Values" comment's
# anchor cell from A4 to A6 along with it.
$commentText = $ws.Range("A4").Comment.Text()

$ws.Rows("4:5").Insert()
$ws.Rows("4:5").RowHeight = 20.1

$ws.Range("A4").Comment.Delete()
$ws.Range("A6").AddComment($commentText)

# Header: rename "type" column header to "color_saturation_type"
$ws.Range("D1").Value = "color_saturation_type"

# Row 2 stays: test_color, 180, 50, low  (unchanged)

# Row 3: was green/120/100/high -> becomes test_color/180/50/medium
$ws.Range("A3").Value = "test_color"
$ws.Range("B3").Value = 180
$ws.Range("C3").Value = 50
$ws.Range("D3").Value = "medium"

# Row 4 (newly inserted): test_color, 180, 50, high
$ws.Range("A4").Value = "test_color"
$ws.Range("B4").Value = 180
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = "high"

# Row 5 (newly inserted): green, 120, 100, high
$ws.Range("A5").Value = "green"
$ws.Range("B5").Value = 120
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = "high"

# Row 6 (old row 4, shifted down): red, 0, 100, high -- values unchanged
